# chore: publish IG 1.0.1
#
# Updates the "Metadata" sheet of the MedCom IHE ClassCode CodeSystem
# spreadsheet for the 1.0.1 IG publication:
#   - Identifier value loses its "id: " prefix
#   - Version bumps from 1.0.0 to 1.0.1
#   - Contact gets a real display value
#   - A new "Jurisdiction" row is inserted right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Identifier: "id: 1.2.208.184.100.9" -> "1.2.208.184.100.9"
$ws.Range("B3").Value = "1.2.208.184.100.9"

# Version: "1.0.0" -> "1.0.1"
$ws.Range("B4").Value = "1.0.1"

# Contact: "No display for ContactDetail" -> "MedCom (http://www.medcom.dk)"
$ws.Range("B11").Value = "MedCom (http://www.medcom.dk)"

# Insert a new "Jurisdiction" row right after "Contact" (row 12), pushing
# "Description" and everything below it down by one row.
$ws.Rows.Item(12).Insert()

# The freshly inserted row doesn't inherit the bordered/wrapped body style
# used throughout the table, so copy the formatting down from the row that
# used to be (and still is, post-insert) row 13.
$ws.Range("A13:B13").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

$ws.Range("A12").Value = "Jurisdiction"

# Force an explicit empty-string text value (rather than a truly blank
# cell) for B12, matching the empty <t/> shared-string entry in the diff.
# A leading apostrophe is Excel's "store as text" prefix; it produces a
# quote-prefixed style, so re-apply the plain body formatting afterwards.
$ws.Range("B12").Value = "'"
$ws.Range("A13:B13").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
